# looks like somebody added plugins
# Update the WA "top remaining" scrape: refresh the scrape date, re-sort a
# few ties that now compare differently, and drop the exhausted
# "$25,000,000 EXTRAVAGANZA" game from the $20.00 Games block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the "2019-03-07" scrape-date column wherever it appears ---
# (keep it a plain text value like the source file, not an auto-converted date)
for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "2019-03-07") {
        $cell.Value = "'2019-03-12"
    }
}

# --- 2. Swap rows 7 & 8: "CASH ON THE RUN" / "DOUBLE IT!" ---
$ws.Cells.Item(7, 3).Value = "DOUBLE IT!"
$ws.Cells.Item(7, 4).Value = 1529
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = "'2019-03-12"

$ws.Cells.Item(8, 3).Value = "CASH ON THE RUN"
$ws.Cells.Item(8, 4).Value = 1462
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = "'2019-02-27"

# --- 3. Swap rows 11 & 12: "MUMMY MAYHEM" / "LUCKY 13" ---
$ws.Cells.Item(11, 3).Value = "LUCKY 13"
$ws.Cells.Item(11, 4).Value = 1494

$ws.Cells.Item(12, 3).Value = "MUMMY MAYHEM"
$ws.Cells.Item(12, 4).Value = 1502

# --- 4. Row 23 "RED BALL TRIPLER": top prizes remaining 3 -> 2 ---
$ws.Cells.Item(23, 5).Value = 2

# --- 5. Re-sort rows 36-42 (the $5.00 Games tie block) into the new order ---
$ws.Cells.Item(36, 3).Value = "MS PAC MAN"
$ws.Cells.Item(36, 4).Value = 1495
$ws.Cells.Item(36, 5).Value = 2

$ws.Cells.Item(37, 3).Value = "GIFTS GALORE"
$ws.Cells.Item(37, 4).Value = 1523
$ws.Cells.Item(37, 5).Value = 1

$ws.Cells.Item(38, 3).Value = "BLACK AND GOLD"
$ws.Cells.Item(38, 4).Value = 1527
$ws.Cells.Item(38, 5).Value = 3

$ws.Cells.Item(39, 3).Value = "SEATTLE SEAHAWKS"
$ws.Cells.Item(39, 4).Value = 1506
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = "'2019-03-12"

$ws.Cells.Item(40, 3).Value = "LADY LUCK"
$ws.Cells.Item(40, 4).Value = 1471
$ws.Cells.Item(40, 5).Value = 1

$ws.Cells.Item(41, 3).Value = "LOTERIA"
$ws.Cells.Item(41, 4).Value = 1515
$ws.Cells.Item(41, 5).Value = 2

$ws.Cells.Item(42, 3).Value = "MUCKLESHOOT CASINO SLOTS"
$ws.Cells.Item(42, 4).Value = 1480
$ws.Cells.Item(42, 5).Value = 3
$ws.Cells.Item(42, 6).Value = "'2019-02-24"

# --- 6. Row 56 "$5 MILLION CASH BLOWOUT": top prizes remaining 613 -> 605 ---
$ws.Cells.Item(56, 5).Value = 605

# --- 7. Drop the exhausted "$25,000,000 EXTRAVAGANZA" row entirely ---
# (it currently sits at row 63; deleting shifts BLACK ICE / 30X CASH up)
$ws.Rows.Item(63).Delete()
